# CLX_cashflow.xlsx update
# - Widen column B (closest achievable to the target stored width of
#   15.400000000000002, matching the width already used by columns C-H).
# - Fill in column B (the most-recent period) with values for every line
#   item that previously had an empty placeholder cell.
# - Update a handful of historical figures in columns C/D/E/F/G for the
#   "Change in inventories" / "Change in payables and accrued liability"
#   rows, plus the restated "Dividends Paid (Common)" figure in B25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width -----------------------------------------------------
$ws.Columns("B").ColumnWidth = 14.65

# --- New period (column B) values -------------------------------------
$ws.Range("B2").Value  = 925000000
$ws.Range("B3").Value  = 533000000
$ws.Range("B4").Value  = -12000000
$ws.Range("B5").Value  = 121000000
$ws.Range("B6").Value  = 234000000
$ws.Range("B7").Value  = 332000000
$ws.Range("B8").Value  = -40000000
$ws.Range("B9").Value  = -72000000
$ws.Range("B10").Value = 1633000000
$ws.Range("B11").Value = -328000000
$ws.Range("B12").Value = -766000000
$ws.Range("B13").Value = -9000000
$ws.Range("B14").Value = -448000000
$ws.Range("B16").Value = -496000000
$ws.Range("B17").Value = -572000000
$ws.Range("B19").Value = -1206000000
$ws.Range("B20").Value = 13000000
$ws.Range("B21").Value = -8000000
$ws.Range("B22").Value = 879000000
$ws.Range("B23").Value = 871000000
$ws.Range("B24").Value = 65000000
$ws.Range("B25").Value = -572000000
$ws.Range("B26").Value = 183000000
$ws.Range("B27").Value = -496000000

# --- Restated historical figures ---------------------------------------
# Change in inventories
$ws.Range("C6").Value = 367000000
$ws.Range("D6").Value = 428000000
$ws.Range("E6").Value = 504000000
$ws.Range("F6").Value = 94000000
$ws.Range("G6").Value = 64000000

# Change in payables and accrued liability
$ws.Range("C7").Value = 394000000
$ws.Range("D7").Value = 393000000
$ws.Range("E7").Value = 291000000
$ws.Range("F7").Value = 125000000
$ws.Range("G7").Value = -20000000
